$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 6, pushing the existing rows 6-34 down to 7-35
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new weekly price record
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C6").Value = "Arica y Parinacota"
$ws.Range("D6").Value = 45114
$ws.Range("E6").Value = 15
$ws.Range("F6").Value = 100112043
$ws.Range("G6").Value = "Pepino dulce"
$ws.Range("H6").Value = "Cultivar XV región"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 160
$ws.Range("K6").Value = 5000
$ws.Range("L6").Value = 6000
$ws.Range("M6").Value = 5500
$ws.Range("N6").Value = "$/caja 10 kilos"
$ws.Range("O6").Value = "Región de Arica y Parinacota"
$ws.Range("P6").Value = 550
$ws.Range("Q6").Value = 10
$ws.Range("R6").Value = "Hortaliza"
